# Generate Report for Handoff
# Regenerate the localization-status report: the handoff file was re-created
# under a new GUID (and a new content hash for the target .xlf files), and
# the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# were refreshed. The hyperlink targets (pointing at a fixed repo commit)
# stay the same - only the visible/display text changes to the new file name.

$wb = $excel.ActiveWorkbook

$oldId = "8ad34f12-178f-48cf-934b-1b7ba2488e68"
$newId = "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee"
$newHash = "ad206ce02866240825132e03ee6a42447766c14c"
$hlAddr = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3899018a2a85bf62f9c9090f8041cc810788f2d9/e2e/" + $oldId + ".md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newId + ".md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hlAddr, [Type]::Missing, [Type]::Missing, "e2e\" + $newId + ".md")

$wsOverview.Range("G2").Value = "2016-08-17 04:51:58"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hlAddr, [Type]::Missing, [Type]::Missing, $newId + ".md")

$wsZhCn.Range("G2").Value = $newId + "." + $newHash + ".zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-17 04:51:54"

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hlAddr, [Type]::Missing, [Type]::Missing, $newId + ".md")

$wsDeDe.Range("G2").Value = $newId + "." + $newHash + ".de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-17 04:51:58"
